# Reorder the "Recorded By" names in column G so that the
# System/system entries move to the end of the comma-separated list
# (instead of leading it), matching the upstream sync of the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, system, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
